# Add three new data columns (China/Italy/UK avg COVID cases) to the monthly
# COVID data sheet, style the new header cells, and restore the selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (F1:H1) ---
$ws.Cells.Item(1, 6).Value = "ChinaAvgCovidCases"
$ws.Cells.Item(1, 7).Value = "ItalyAvgCovidCases"
$ws.Cells.Item(1, 8).Value = "UKAvgCovidCases"

# --- Data rows 2-40 for the new columns ---
$ws.Cells.Item(2, 6).Value = 335.17241379310337
$ws.Cells.Item(2, 7).Value = 0.125
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(3, 6).Value = 2402.379310344827
$ws.Cells.Item(3, 7).Value = 30.517241379310349
$ws.Cells.Item(3, 8).Value = 1.931034482758621
$ws.Cells.Item(4, 6).Value = 101.8064516129032
$ws.Cells.Item(4, 7).Value = 3253.2580645161288
$ws.Cells.Item(4, 8).Value = 945.70967741935488
$ws.Cells.Item(5, 6).Value = 60.93333333333333
$ws.Cells.Item(5, 7).Value = 3395.0666666666671
$ws.Cells.Item(5, 8).Value = 4567.3999999999996
$ws.Cells.Item(6, 6).Value = 6.354838709677419
$ws.Cells.Item(6, 7).Value = 937.83870967741939
$ws.Cells.Item(6, 8).Value = 2833.0645161290322
$ws.Cells.Item(7, 6).Value = 21.9
$ws.Cells.Item(7, 7).Value = 272.36666666666667
$ws.Cells.Item(7, 8).Value = 973.33333333333337
$ws.Cells.Item(8, 6).Value = 88.032258064516128
$ws.Cells.Item(8, 7).Value = 216.83870967741939
$ws.Cells.Item(8, 8).Value = 665.9677419354839
$ws.Cells.Item(9, 6).Value = 78.290322580645167
$ws.Cells.Item(9, 7).Value = 679.35483870967744
$ws.Cells.Item(9, 8).Value = 1086.161290322581
$ws.Cells.Item(10, 6).Value = 21.93333333333333
$ws.Cells.Item(10, 7).Value = 1493.1
$ws.Cells.Item(10, 8).Value = 4366.2
$ws.Cells.Item(11, 6).Value = 27.483870967741939
$ws.Cells.Item(11, 7).Value = 10795.58064516129
$ws.Cells.Item(11, 8).Value = 18175.77419354839
$ws.Cells.Item(12, 6).Value = 52.4
$ws.Cells.Item(12, 7).Value = 31250.133333333339
$ws.Cells.Item(12, 8).Value = 20629.633333333339
$ws.Cells.Item(13, 6).Value = 103.48387096774189
$ws.Cells.Item(13, 7).Value = 16068.129032258061
$ws.Cells.Item(13, 8).Value = 29435.516129032261
$ws.Cells.Item(14, 6).Value = 138.741935483871
$ws.Cells.Item(14, 7).Value = 14777.22580645161
$ws.Cells.Item(14, 8).Value = 42489.129032258068
$ws.Cells.Item(15, 6).Value = 33.785714285714278
$ws.Cells.Item(15, 7).Value = 13072.928571428571
$ws.Cells.Item(15, 8).Value = 13091.642857142861
$ws.Cells.Item(16, 6).Value = 27.161290322580641
$ws.Cells.Item(16, 7).Value = 21068.38709677419
$ws.Cells.Item(16, 8).Value = 5513.8709677419356
$ws.Cells.Item(17, 6).Value = 28.666666666666671
$ws.Cells.Item(17, 7).Value = 14942.1
$ws.Cells.Item(17, 8).Value = 2651.333333333333
$ws.Cells.Item(18, 6).Value = 254.87096774193549
$ws.Cells.Item(18, 7).Value = 6670.8064516129034
$ws.Cells.Item(18, 8).Value = 2462.2258064516132
$ws.Cells.Item(19, 6).Value = 234.9666666666667
$ws.Cells.Item(19, 7).Value = 1437.666666666667
$ws.Cells.Item(19, 8).Value = 10617.2
$ws.Cells.Item(20, 6).Value = 69.354838709677423
$ws.Cells.Item(20, 7).Value = 2722.1290322580639
$ws.Cells.Item(20, 8).Value = 34886.483870967742
$ws.Cells.Item(21, 6).Value = 75.677419354838705
$ws.Cells.Item(21, 7).Value = 6160.6451612903229
$ws.Cells.Item(21, 8).Value = 30743.129032258061
$ws.Cells.Item(22, 6).Value = 53.5
$ws.Cells.Item(22, 7).Value = 4458.7333333333336
$ws.Cells.Item(22, 8).Value = 34492.9
$ws.Cells.Item(23, 6).Value = 45.322580645161288
$ws.Cells.Item(23, 7).Value = 3199.322580645161
$ws.Cells.Item(23, 8).Value = 40996.548387096773
$ws.Cells.Item(24, 6).Value = 62
$ws.Cells.Item(24, 7).Value = 8278.3333333333339
$ws.Cells.Item(24, 8).Value = 39566.133333333331
$ws.Cells.Item(25, 6).Value = 133.32258064516131
$ws.Cells.Item(25, 7).Value = 31149.61290322581
$ws.Cells.Item(25, 8).Value = 97481.193548387091
$ws.Cells.Item(26, 6).Value = 230.12903225806451
$ws.Cells.Item(26, 7).Value = 159485.70967741939
$ws.Cells.Item(26, 8).Value = 128591.6451612903
$ws.Cells.Item(27, 6).Value = 7023.9642857142853
$ws.Cells.Item(27, 7).Value = 65681.178571428565
$ws.Cells.Item(27, 8).Value = 54571.357142857138
$ws.Cells.Item(28, 6).Value = 18511.193548387098
$ws.Cells.Item(28, 7).Value = 58175.225806451614
$ws.Cells.Item(28, 8).Value = 70323.612903225803
$ws.Cells.Item(29, 6).Value = 6045.2333333333336
$ws.Cells.Item(29, 7).Value = 61373.1
$ws.Cells.Item(29, 8).Value = 30526.73333333333
$ws.Cells.Item(30, 6).Value = 62202.032258064522
$ws.Cells.Item(30, 7).Value = 31856.129032258061
$ws.Cells.Item(30, 8).Value = 8032.7419354838712
$ws.Cells.Item(31, 6).Value = 58334.133333333331
$ws.Cells.Item(31, 7).Value = 34738.466666666667
$ws.Cells.Item(31, 8).Value = 14085.73333333333
$ws.Cells.Item(32, 6).Value = 27314.61290322581
$ws.Cells.Item(32, 7).Value = 82706.322580645166
$ws.Cells.Item(32, 8).Value = 19983.967741935481
$ws.Cells.Item(33, 6).Value = 25334
$ws.Cells.Item(33, 7).Value = 27199.032258064519
$ws.Cells.Item(33, 8).Value = 5762.2258064516127
$ws.Cells.Item(34, 6).Value = 40566.73333333333
$ws.Cells.Item(34, 7).Value = 19562
$ws.Cells.Item(34, 8).Value = 5323.3
$ws.Cells.Item(35, 6).Value = 42870.93548387097
$ws.Cells.Item(35, 7).Value = 36023.225806451614
$ws.Cells.Item(35, 8).Value = 7665.322580645161
$ws.Cells.Item(36, 6).Value = 23973.3
$ws.Cells.Item(36, 7).Value = 28832.76666666667
$ws.Cells.Item(36, 8).Value = 3464.1
$ws.Cells.Item(37, 6).Value = 2427689.2903225808
$ws.Cells.Item(37, 7).Value = 24320.38709677419
$ws.Cells.Item(37, 8).Value = 5538.0967741935483
$ws.Cells.Item(38, 6).Value = 437741.83870967739
$ws.Cells.Item(38, 7).Value = 9730.0645161290322
$ws.Cells.Item(38, 8).Value = 3035.1290322580639
$ws.Cells.Item(39, 6).Value = 19110.321428571431
$ws.Cells.Item(39, 7).Value = 4265.1071428571431
$ws.Cells.Item(39, 8).Value = 3631.821428571428
$ws.Cells.Item(40, 6).Value = 6737.0967741935483
$ws.Cells.Item(40, 7).Value = 3413.483870967742
$ws.Cells.Item(40, 8).Value = 4341.6129032258068

# --- Header style: bold, centered/top aligned, boxed border ---
$hdr = $ws.Range("F1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.LineStyle = 1         # xlContinuous (thin box)

# --- Row 1 height / new column widths ---
$ws.Rows.Item(1).RowHeight = 38.4
$ws.Columns.Item(6).AutoFit()
$ws.Columns.Item(7).AutoFit()
$ws.Columns.Item(8).AutoFit()

# --- Restore final selection ---
$ws.Range("H12").Select()

Write-Host "Done"
